# "Changed to use raw access to hardware / improved motion algorithm for
# stepper motors" - update the trapezoid-planner inputs on Sheet1. All of
# the downstream formulas (accDist, decelDist, plateauDist, the
# accel/plateau/finish points used by the chart, etc.) recalculate
# automatically from these five named inputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 50    # feedrate
$ws.Range("C5").Value = 0     # init
$ws.Range("C6").Value = 0     # final
$ws.Range("C7").Value = 100   # acc
$ws.Range("C8").Value = 100   # dist

# Recalculate so every formula cell (and the chart that plots C25:D28)
# picks up the new inputs.
$excel.CalculateFullRebuild()

# Reproduce the recorded UI state: the active cell/selection moved to C8.
$ws.Range("C8").Select() | Out-Null

# Best-effort: reproduce the recorded application window size.
$excel.Width = 23445
$excel.Height = 12420
$win = $wb.Windows.Item(1)
$win.Width = 23445
$win.Height = 12420
